$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.273.43'
$ws.Range("E2").Value = '  -2.62%  '

$ws.Range("D3").Value = '3.174.29'
$ws.Range("E3").Value = '  -7.94%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.77'
$ws.Range("E5").Value = '  -3.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.92'
$ws.Range("E6").Value = '  -1.71%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  +1.00%  '

$ws.Range("D9").Value = '3.176.18'
$ws.Range("E9").Value = '  -7.84%  '

$ws.Range("E10").Value = '  -6.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.61'
$ws.Range("E11").Value = '  -4.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.395'
$ws.Range("E12").Value = '  -4.33%  '

$ws.Range("D13").Value = '3.726.24'
$ws.Range("E13").Value = '  -7.86%  '

$ws.Range("E14").Value = '  +1.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.38'
$ws.Range("E15").Value = '  -6.13%  '

$ws.Range("D16").Value = '64.286.86'
$ws.Range("E16").Value = '  -2.56%  '

$ws.Range("E17").Value = '  -5.06%  '

$ws.Range("D18").Value = '3.171.05'
$ws.Range("E18").Value = '  -7.89%  '

$ws.Range("E19").Value = '  -4.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.03'
$ws.Range("E20").Value = '  -5.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.91'
$ws.Range("E21").Value = '  -4.75%  '

$ws.Range("E22").Value = '  -5.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.17'
$ws.Range("E24").Value = '  -4.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.503'
$ws.Range("E25").Value = '  -5.26%  '

$ws.Range("E26").Value = '  -3.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.63'
$ws.Range("E27").Value = '  -1.06%  '

$ws.Range("E28").Value = '  -1.51%  '

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.63'
$ws.Range("E30").Value = '  -3.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.11%  '

$ws.Range("E32").Value = '  -4.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.12'
$ws.Range("E33").Value = '  -6.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.66'
$ws.Range("E34").Value = '  -5.01%  '

$ws.Range("E35").Value = '  -6.02%  '

$ws.Range("E36").Value = '  -5.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.81'
$ws.Range("E37").Value = '  -4.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.811'
$ws.Range("E38").Value = '  -7.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.07'
$ws.Range("E39").Value = '  -8.46%  '

$ws.Range("E40").Value = '  -2.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.68'
$ws.Range("E41").Value = '  -6.17%  '

$ws.Range("D42").Value = '2.630.78'
$ws.Range("E42").Value = '  -5.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.17'
$ws.Range("E43").Value = '  -6.70%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '334.70'
$ws.Range("E44").Value = '  +2.02%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.01'
$ws.Range("E45").Value = '  -6.91%  '

$ws.Range("E46").Value = '  -4.29%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '38.88'
$ws.Range("E47").Value = '  -2.92%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.90'
$ws.Range("E48").Value = '  -4.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0270'
$ws.Range("E49").Value = '  -7.26%  '

$ws.Range("E50").Value = '  -1.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.03%  '
